$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I5").Value = 96.21122363135699
$ws.Range("I6").Value = 96.28043801996
$ws.Range("I7").Value = 96.34734763760507

$ws.Range("G20").Value = 98.05481488540376
$ws.Range("G21").Value = 98.14687494156553
$ws.Range("G22").Value = 98.0560433091354

$ws.Range("H23").Value = 97.59017290437959
$ws.Range("H24").Value = 97.539849727918
$ws.Range("H25").Value = 97.58252248064257

$ws.Range("I28").Value = 96.0959689343416
$ws.Range("I29").Value = 96.10059567384077

$ws.Range("G38").Value = 97.97602838626219
$ws.Range("G39").Value = 98.01453679926352

$ws.Range("H40").Value = 97.67522593120876
$ws.Range("H41").Value = 97.80876634216936

$ws.Range("I44").Value = 96.29339976020957
$ws.Range("I45").Value = 96.41247918002682

$ws.Range("G54").Value = 98.01557500493178
$ws.Range("G55").Value = 98.05135217399311

$ws.Range("H56").Value = 97.49088800122603
$ws.Range("H57").Value = 97.53581224532265
